$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.302.41'
$ws.Range('E2').Value = '  +1.36%  '
$ws.Range('D3').Value = '2.636.70'
$ws.Range('E3').Value = '  +1.01%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '599.40'
$ws.Range('E5').Value = '  +1.26%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '154.18'
$ws.Range('E6').Value = '  +2.33%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.545'
$ws.Range('E8').Value = '  -0.08%  '
$ws.Range('D9').Value = '2.635.57'
$ws.Range('E9').Value = '  +0.99%  '
$ws.Range('E10').Value = '  +7.75%  '
$ws.Range('E11').Value = '  -0.66%  '
$ws.Range('E12').Value = '  +0.94%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.348'
$ws.Range('E13').Value = '  +1.47%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.92'
$ws.Range('E14').Value = '  +2.34%  '
$ws.Range('E15').Value = '  +2.79%  '
$ws.Range('D16').Value = '3.128.36'
$ws.Range('E16').Value = '  +1.52%  '
$ws.Range('D17').Value = '68.277.17'
$ws.Range('E17').Value = '  +1.65%  '
$ws.Range('D18').Value = '2.635.18'
$ws.Range('E18').Value = '  +1.01%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.39'
$ws.Range('E19').Value = '  +3.18%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '365.85'
$ws.Range('E20').Value = '  -1.62%  '
$ws.Range('E22').Value = '  -0.64%  '
$ws.Range('E23').Value = '  +1.76%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.11'
$ws.Range('E24').Value = '  +3.46%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '73.67'
$ws.Range('E25').Value = '  -0.01%  '
$ws.Range('E26').Value = '  -0.02%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.03'
$ws.Range('E27').Value = '  +1.20%  '
$ws.Range('D28').Value = '2.770.29'
$ws.Range('E28').Value = '  +1.00%  '
$ws.Range('E29').Value = '  +5.67%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.999'
$ws.Range('E30').Value = '  -0.17%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '573.28'
$ws.Range('E31').Value = '  -1.23%  '
$ws.Range('E32').Value = '  +4.02%  '
$ws.Range('E33').Value = '  +4.16%  '
$ws.Range('E34').Value = '  +2.35%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.129'
$ws.Range('E35').Value = '  +2.88%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.999'
$ws.Range('E36').Value = '  +0.01%  '
$ws.Range('E37').Value = '  +3.28%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '160.30'
$ws.Range('E38').Value = '  +1.64%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '19.24'
$ws.Range('E39').Value = '  +0.95%  '
$ws.Range('E40').Value = '  +3.19%  '
$ws.Range('E41').Value = '  +0.66%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.37'
$ws.Range('E42').Value = '  +2.64%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '17.73'
$ws.Range('E43').Value = '  +3.46%  '
$ws.Range('E44').Value = '  +2.70%  '
$ws.Range('D45').Value = '0.0₆0320'
$ws.Range('E45').Value = '  +12.91%  '
$ws.Range('E46').Value = '  +0.04%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '40.49'
$ws.Range('E47').Value = '  -0.42%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '157.07'
$ws.Range('E48').Value = '  +2.51%  '
$ws.Range('E49').Value = '  +0.75%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.70'
$ws.Range('E50').Value = '  +1.75%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '21.85'
$ws.Range('E51').Value = '  +2.21%  '
